$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts existing B,C,D -> C,D,E)
$ws.Range("B1").EntireColumn.Insert()

# Set header
$ws.Range("B1").Value = "range_end"

# Fill in range_end values: pattern 30,50,70,100 repeating for rows 2-53
$values = @("30", "50", "70", "100")
$rangeEndCol = $ws.Range("B2:B53")
$rangeEndCol.NumberFormat = "@"
for ($row = 2; $row -le 53; $row++) {
    $idx = ($row - 2) % 4
    $ws.Cells.Item($row, 2).Value = $values[$idx]
}
